$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
# A8 shared string run "3" -> "4" (Volume 33   Number  4)
$ws.Range("A8").Characters(21, 1).Text = "4"
# C9 shared string runs: week-of dates shift by one week
$ws.Range("C9").Characters(27, 9).Text = "1/19/2026"
$ws.Range("C9").Characters(47, 9).Text = "1/25/2026"

# --- Crime statistics table (rows 14-31) ---
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "***.*"

# Row 15
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 66.666666666666
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 66.666666666666
$ws.Range("L15").Value = 66.666666666666
$ws.Range("M15").Value = 66.666666666666

# Row 16
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 13
$ws.Range("E16").Value = -46.153846153846
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 36
$ws.Range("H16").Value = -44.444444444444
$ws.Range("I16").Value = 20
$ws.Range("J16").Value = 35
$ws.Range("K16").Value = -42.857142857142
$ws.Range("L16").Value = -31.034482758620
$ws.Range("M16").Value = -13.043478260869
$ws.Range("N16").Value = -87.012987012987

# Row 17
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 55
$ws.Range("G17").Value = 45
$ws.Range("H17").Value = 22.222222222222
$ws.Range("I17").Value = 52
$ws.Range("J17").Value = 35
$ws.Range("K17").Value = 48.571428571428
$ws.Range("L17").Value = 36.842105263157
$ws.Range("M17").Value = 92.592592592592
$ws.Range("N17").Value = -13.333333333333

# Row 18
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 31
$ws.Range("H18").Value = -58.064516129032
$ws.Range("I18").Value = 9
$ws.Range("J18").Value = 29
$ws.Range("K18").Value = -68.965517241379
$ws.Range("L18").Value = -43.75
$ws.Range("M18").Value = -74.285714285714
$ws.Range("N18").Value = -91.346153846153

# Row 19
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 28
$ws.Range("E19").Value = -46.428571428571
$ws.Range("F19").Value = 64
$ws.Range("G19").Value = 72
$ws.Range("H19").Value = -11.111111111111
$ws.Range("I19").Value = 56
$ws.Range("J19").Value = 66
$ws.Range("K19").Value = -15.151515151515
$ws.Range("L19").Value = -11.111111111111
$ws.Range("M19").Value = 80.645161290322
$ws.Range("N19").Value = 47.368421052631

# Row 20
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 22.222222222222
$ws.Range("F20").Value = 34
$ws.Range("H20").Value = -15
$ws.Range("I20").Value = 30
$ws.Range("J20").Value = 32
$ws.Range("K20").Value = -6.25
$ws.Range("L20").Value = -11.764705882352
$ws.Range("M20").Value = 130.769230769231
$ws.Range("N20").Value = -81.927710843373

# Row 21
$ws.Range("C21").Value = 52
$ws.Range("D21").Value = 72
$ws.Range("E21").Value = -27.777777777777
$ws.Range("F21").Value = 191
$ws.Range("G21").Value = 229
$ws.Range("H21").Value = -16.593886462882
$ws.Range("I21").Value = 172
$ws.Range("J21").Value = 202
$ws.Range("K21").Value = -14.851485148514
$ws.Range("L21").Value = -6.010928961748
$ws.Range("M21").Value = 30.303030303030
$ws.Range("N21").Value = -67.362428842504

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("G22").Value = 2

# Row 23
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 16
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = 33.333333333333
$ws.Range("I23").Value = 16
$ws.Range("J23").Value = 9
$ws.Range("K23").Value = 77.777777777777
$ws.Range("L23").Value = -20
$ws.Range("M23").Value = 33.333333333333

# Row 24
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 41
$ws.Range("E24").Value = -26.829268292682
$ws.Range("F24").Value = 144
$ws.Range("G24").Value = 161
$ws.Range("H24").Value = -10.559006211180
$ws.Range("I24").Value = 129
$ws.Range("J24").Value = 137
$ws.Range("K24").Value = -5.839416058394
$ws.Range("L24").Value = 2.380952380952
$ws.Range("M24").Value = 30.303030303030

# Row 25
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -31.25
$ws.Range("G25").Value = 47
$ws.Range("H25").Value = -27.659574468085
$ws.Range("I25").Value = 26
$ws.Range("J25").Value = 41
$ws.Range("K25").Value = -36.585365853658
$ws.Range("L25").Value = -43.478260869565

# Row 26
$ws.Range("C26").Value = 21
$ws.Range("D26").Value = 18
$ws.Range("E26").Value = 16.666666666666
$ws.Range("F26").Value = 83
$ws.Range("G26").Value = 82
$ws.Range("H26").Value = 1.219512195121
$ws.Range("I26").Value = 79
$ws.Range("J26").Value = 66
$ws.Range("K26").Value = 19.696969696969
$ws.Range("L26").Value = 14.492753623188
$ws.Range("M26").Value = -14.130434782608

# Row 27
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 100
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = 100
$ws.Range("L27").Value = 50

# Row 28
$ws.Range("C28").Value = 3
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Value = 2
$ws.Range("E28").NumberFormat = "General"
$ws.Range("E28").Value = 50
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 25
$ws.Range("I28").Value = 4
$ws.Range("J28").Value = 3
$ws.Range("K28").Value = 33.333333333333
$ws.Range("L28").Value = -60

# Row 29
$ws.Range("J29").Value = 2

# Row 30
$ws.Range("J30").Value = 2

# Row 31
$ws.Range("F31").NumberFormat = "@"
$ws.Range("F31").Value = "0"

